$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 25 (shifts old rows 25-48 down to 26-49),
# then populate it with the new weekly price observation.
$ws.Rows.Item(25).Insert()

$ws.Cells.Item(25,1).Value2 = 11
$ws.Cells.Item(25,2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(25,3).Value2 = "Bíobío"
$ws.Cells.Item(25,4).Value2 = 45128
$ws.Cells.Item(25,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(25,5).Value2 = 8
$ws.Cells.Item(25,6).Value2 = "Fruta"
$ws.Cells.Item(25,7).Value2 = 100107
$ws.Cells.Item(25,8).Value2 = "Otros"
$ws.Cells.Item(25,9).Value2 = 100107002
$ws.Cells.Item(25,10).Value2 = "Chirimoya"
$ws.Cells.Item(25,11).Value2 = "Cultivar IV Región"
$ws.Cells.Item(25,12).Value2 = "Primera"
$ws.Cells.Item(25,13).Value2 = 50
$ws.Cells.Item(25,14).Value2 = 30000
$ws.Cells.Item(25,15).Value2 = 30000
$ws.Cells.Item(25,16).Value2 = 30000
$ws.Cells.Item(25,17).Value2 = "`$/bandeja 10 kilos"
$ws.Cells.Item(25,18).Value2 = "Provincia de Limarí"
$ws.Cells.Item(25,19).Value2 = 3000
$ws.Cells.Item(25,20).Value2 = 10
